$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (D), new Volume(1h) (E), new Hora (G).
# An empty string means "leave this column unchanged" (price/volume
# are blank for the placeholder "--" / "--%" rows and for row 8 price).
$updates = @(
    @{Row=2; D="322.43"; E="-2.09%"; G="8"},
    @{Row=3; D="39.38"; E="-1.57%"; G="8"},
    @{Row=4; D="5.711"; E="8.65%"; G="8"},
    @{Row=5; D="0.07999"; E="-1.19%"; G="8"},
    @{Row=6; D="8.614"; E="-0.36%"; G="8"},
    @{Row=7; D="1.961"; E="1.88%"; G="8"},
    @{Row=8; D=""; E="-0.25%"; G="8"},
    @{Row=9; D="0.9272"; E="-0.91%"; G="8"},
    @{Row=10; D="0.1272"; E="-4.29%"; G="8"},
    @{Row=11; D="0.1945"; E="-1.48%"; G="8"},
    @{Row=12; D="8.735"; E="25.28%"; G="8"},
    @{Row=13; D="0.09128"; E="0.59%"; G="8"},
    @{Row=14; D="0.03679"; E="4.03%"; G="8"},
    @{Row=15; D="0.1049"; E="9.50%"; G="8"},
    @{Row=16; D="0.001294"; E="-3.45%"; G="8"},
    @{Row=17; D="0.006355"; E="5.34%"; G="8"},
    @{Row=18; D="3.351"; E="-0.51%"; G="8"},
    @{Row=19; D="4.550"; E="0.34%"; G="8"},
    @{Row=20; D="0.3538"; E="0.84%"; G="8"},
    @{Row=21; D="0.1374"; E="3.95%"; G="8"},
    @{Row=22; D="0.2452"; E="-4.62%"; G="8"},
    @{Row=23; D="0.04414"; E="-0.73%"; G="8"},
    @{Row=24; D="0.001264"; E="3.38%"; G="8"},
    @{Row=25; D="0.004526"; E="5.02%"; G="8"},
    @{Row=26; D="0.0001153"; E="-3.23%"; G="8"},
    @{Row=27; D=""; E=""; G="8"},
    @{Row=28; D=""; E=""; G="8"},
    @{Row=29; D=""; E=""; G="8"},
    @{Row=30; D=""; E=""; G="8"},
    @{Row=31; D=""; E=""; G="8"},
    @{Row=32; D=""; E=""; G="8"},
    @{Row=33; D=""; E=""; G="8"},
    @{Row=34; D=""; E=""; G="8"},
    @{Row=35; D=""; E=""; G="8"},
    @{Row=36; D=""; E=""; G="8"},
    @{Row=37; D=""; E=""; G="8"},
    @{Row=38; D=""; E=""; G="8"},
    @{Row=39; D="0.02493"; E="-0.47%"; G="8"},
    @{Row=40; D="0.05362"; E="3.46%"; G="8"},
    @{Row=41; D="0.007460"; E="-3.68%"; G="8"},
    @{Row=42; D="0.009574"; E="4.73%"; G="8"},
    @{Row=43; D="0.1402"; E="-1.63%"; G="8"},
    @{Row=44; D="0.002121"; E="-1.91%"; G="8"},
    @{Row=45; D="0.01096"; E="-1.58%"; G="8"},
    @{Row=46; D="0.00006761"; E="1.56%"; G="8"},
    @{Row=47; D="0.00000000752"; E="0.13%"; G="8"},
    @{Row=48; D="0.002978"; E="-10.96%"; G="8"},
    @{Row=49; D="0.002296"; E="-7.53%"; G="8"},
    @{Row=50; D="0.00002105"; E="0.13%"; G="8"},
    @{Row=51; D="0.0002004"; E="0.13%"; G="8"}
)

foreach ($u in $updates) {
    foreach ($col in @("D", "E", "G")) {
        $newVal = $u[$col]
        if ([string]::IsNullOrEmpty($newVal)) { continue }
        $cell = $ws.Range("$col$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $newVal
    }
}
